# CI: Auto Update Data (#88)
# Applies the upstream data refresh to the "悖论模拟干员名单用户版" sheet:
#  - updates several "count"/"codes" cell pairs with revised maa:// links
#  - adds a handful of brand-new operator rows/entries
#  - fills in previously blank "0 / None" cells

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $value) {
    # Force literal text storage (mirrors how these sheets store every cell,
    # including purely-numeric-looking counts, as text) so values such as
    # "2" or "0" don't silently become numbers.
    $cell = $ws.Range($rangeAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- AC/AD/AE (特种) row 2: 砾 ---
Set-TextValue "AD2" "2"
$ws.Range("AE2").Value = "maa://25251, *maa://36675"

# --- U/V/W (医疗) row 4: 芙蓉 ---
Set-TextValue "V4" "3"
$ws.Range("W4").Value = "**maa://32495, ***maa://31785, ***maa://36683"

# --- I/J/K (重装) row 5: 卡缇 ---
$ws.Range("K5").Value = "*maa://22757"

# --- Q/R/S (术师) row 6: 12F ---
$ws.Range("S6").Value = "maa://37411"

# --- AC/AD/AE (特种) row 7 ---
$ws.Range("AE7").Value = "*maa://26191, *maa://36671, *maa://42530"

# --- I/J/K (重装) row 9 ---
$ws.Range("K9").Value = "maa://22762, maa://39552"

# --- Q/R/S (术师) row 16 ---
$ws.Range("S16").Value = "maa://22729, *maa://28648, maa://36674"

# --- Q/R/S (术师) row 17 ---
$ws.Range("S17").Value = "**maa://42324"

# --- A/B/C (先锋) row 20 ---
$ws.Range("C20").Value = "maa://21432, maa://25198, *maa://20795, maa://36680"

# --- Y/Z/AA (辅助) row 21 ---
$ws.Range("AA21").Value = "*maa://21443, ***maa://23820"

# --- U/V/W (医疗) row 24 ---
$ws.Range("W24").Value = "maa://29988, maa://23504, **maa://22892, *maa://25141, *maa://36663, ***maa://22815"

# --- E/F/G (近卫) row 27 ---
$ws.Range("G27").Value = "**maa://21283, maa://34494, *maa://39601, **maa://36665"

# --- AC/AD/AE (特种) row 29 ---
Set-TextValue "AD29" "2"
$ws.Range("AE29").Value = "*maa://24080, maa://42865"

# --- Y/Z/AA (辅助) row 30 ---
Set-TextValue "Z30" "1"
$ws.Range("AA30").Value = "maa://42979"

# --- I/J/K (重装) row 31 ---
$ws.Range("K31").Value = "maa://35926, *maa://36258"

# --- Q/R/S (术师) row 32 ---
Set-TextValue "R32" "3"
$ws.Range("S32").Value = "maa://41108, maa://41238, maa://42859"

# --- New operator: 忍冬 (先锋) on row 35, no existing A/B/C entry there yet ---
$ws.Range("A35").Value = "忍冬"
$ws.Range("B35").Value = "-"
$ws.Range("C35").Value = "-"

# --- New operator: 云迹 (特种) on row 40 ---
$ws.Range("AC40").Value = "云迹"
$ws.Range("AD40").Value = "-"
$ws.Range("AE40").Value = "-"

# --- New operator: 菲莱 (重装) on row 41 ---
$ws.Range("I41").Value = "菲莱"
$ws.Range("J41").Value = "-"
$ws.Range("K41").Value = "-"

# --- New operator: 裁度 (特种) on row 41 ---
$ws.Range("AC41").Value = "裁度"
$ws.Range("AD41").Value = "-"
$ws.Range("AE41").Value = "-"

# --- New operator: 弑君者 (特种) on row 42 ---
$ws.Range("AC42").Value = "弑君者"
$ws.Range("AD42").Value = "-"
$ws.Range("AE42").Value = "-"

# --- M/N/O (狙击) row 43: 截云 now has data (0 / None) ---
Set-TextValue "N43" "0"
$ws.Range("O43").Value = "None"

# --- Q/R/S (术师) row 43: 折光 now has data (0 / None) ---
Set-TextValue "R43" "0"
$ws.Range("S43").Value = "None"

# --- New operator: 荒芜拉普兰德 (术师) on row 49 ---
$ws.Range("Q49").Value = "荒芜拉普兰德"
$ws.Range("R49").Value = "-"
$ws.Range("S49").Value = "-"

# --- E/F/G (近卫) row 62: 止颂 now has data ---
Set-TextValue "F62" "1"
$ws.Range("G62").Value = "maa://42981"
